$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.909.30'
$ws.Range('E2').Value = '  -2.05%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.993.50'
$ws.Range('E3').Value = '  -3.11%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.35'
$ws.Range('E5').Value = '  -0.40%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.640'
$ws.Range('E6').Value = '  -3.86%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.53'
$ws.Range('E7').Value = '  +8.30%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.87'
$ws.Range('E9').Value = '  -3.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.366'
$ws.Range('E10').Value = '  -0.45%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0742'
$ws.Range('E11').Value = '  -1.83%  '

$ws.Range('E12').Value = '  -2.10%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.946'
$ws.Range('E13').Value = '  -3.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.80'
$ws.Range('E14').Value = '  -0.24%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.282.40'
$ws.Range('E15').Value = '  -3.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.35'
$ws.Range('E16').Value = '  -2.43%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.66'
$ws.Range('E17').Value = '  +12.86%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.007.75'
$ws.Range('E18').Value = '  -2.46%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '35.809.93'
$ws.Range('E19').Value = '  -2.13%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.98'
$ws.Range('E20').Value = '  -0.43%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0851'
$ws.Range('E21').Value = '  -1.37%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.23'
$ws.Range('E22').Value = '  -0.74%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.68'
$ws.Range('E23').Value = '  -2.03%  '

$ws.Range('E24').Value = '  +0.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.62'
$ws.Range('E25').Value = '  +12.85%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('E26').Value = '  -4.34%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.68'
$ws.Range('E27').Value = '  +4.27%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.38'
$ws.Range('E28').Value = '  -0.63%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.41'
$ws.Range('E29').Value = '  -3.79%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.120'
$ws.Range('E30').Value = '  -2.01%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.95'
$ws.Range('E31').Value = '  -2.85%  '

$ws.Range('E32').Value = '  -6.39%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0984'
$ws.Range('E33').Value = '  +14.09%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0606'
$ws.Range('E34').Value = '  +1.47%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').Value = '  +9.27%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.42'
$ws.Range('E36').Value = '  -2.15%  '

$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('E38').Value = '  -1.27%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.79'
$ws.Range('E39').Value = '  +13.99%  '

$ws.Range('E40').Value = '  -1.46%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0214'
$ws.Range('E41').Value = '  -0.93%  '

$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.86'
$ws.Range('E42').Value = '  -1.44%  '

$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0938'
$ws.Range('E43').Value = '  +2.08%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.11'
$ws.Range('E44').Value = '  -0.26%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '94.42'
$ws.Range('E45').Value = '  -0.81%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.52'
$ws.Range('E46').Value = '  +3.20%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.82'
$ws.Range('E47').Value = '  +2.82%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.373.66'
$ws.Range('E48').Value = '  -2.90%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.34'
$ws.Range('E49').Value = '  +2.77%  '

$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.90'
$ws.Range('E50').Value = '  -1.48%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.89'
$ws.Range('E51').Value = '  +1.86%  '
